$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.831.99"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.089.23"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.29%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.25%  "
$ws.Range("E6").Value = "  -0.25%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.47"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.61%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +0.66%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0784"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.66%  "
$ws.Range("E11").Value = "  +2.70%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.26"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.91%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.397.95"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.32%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.26"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.57%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.781"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.78%  "
$ws.Range("E16").Value = "  +1.06%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.099.72"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.84%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.829.32"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.33%  "
$ws.Range("E19").Value = "  -0.26%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.06"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0838"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "229.97"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.44%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("E25").Value = "  -1.25%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.77"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +8.36%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "171.71"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.41%  "
$ws.Range("E28").Value = "  -3.15%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.54"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.41"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.11%  "
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("E32").Value = "  +0.08%  "
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.63"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.90%  "
$ws.Range("E35").Value = "  +0.26%  "
$ws.Range("E36").Value = "  -0.81%  "
$ws.Range("E37").Value = "  -1.77%  "
$ws.Range("E38").Value = "  -0.10%  "
$ws.Range("E39").Value = "  -0.37%  "
$ws.Range("E40").Value = "  +8.63%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "101.36"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.80%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0974"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.96%  "
$ws.Range("E43").Value = "  +1.87%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.19"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.50%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.80"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.453.84"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.62%  "
$ws.Range("E47").Value = "  -0.47%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.11"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.67%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.21"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.81%  "
$ws.Range("E50").Value = "  -1.80%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.281.45"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.26%  "
